$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'20.392.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.42%  "

$ws.Range("D3").Value = "'1.440.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.51%  "

$ws.Range("E4").Value = "  -0.40%  "

$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").Value = "'277.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.85%  "

$ws.Range("D7").Value = "'0.3730"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.82%  "

$ws.Range("D8").Value = "'0.3091"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.13%  "

$ws.Range("D9").Value = "'40.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.30%  "

$ws.Range("D10").Value = "'1.014"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.48%  "

$ws.Range("D11").Value = "'0.06598"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.76%  "

$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").Value = "'5.376"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.41%  "

$ws.Range("D14").Value = "'17.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.51%  "

$ws.Range("D15").Value = "'6.150"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.23%  "

$ws.Range("D16").Value = "'1.438.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.19%  "

$ws.Range("D17").Value = "'0.00001014"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.41%  "

$ws.Range("D18").Value = "'76.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.84%  "

$ws.Range("D19").Value = "'0.05836"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -11.13%  "

$ws.Range("D21").Value = "'5.746"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.33%  "

$ws.Range("E22").Value = "  -5.37%  "

$ws.Range("D23").Value = "'10.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.51%  "

$ws.Range("D24").Value = "'2.321"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.04%  "

$ws.Range("D25").Value = "'20.383.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.54%  "

$ws.Range("D26").Value = "'2.251"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.98%  "

$ws.Range("D27").Value = "'142.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.48%  "

$ws.Range("D28").Value = "'17.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.04%  "

$ws.Range("D29").Value = "'1.603.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.99%  "

$ws.Range("D30").Value = "'110.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.62%  "

$ws.Range("D31").Value = "'3.965"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -18.31%  "

$ws.Range("D32").Value = "'0.9230"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.38%  "

$ws.Range("D33").Value = "'5.509"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.75%  "

$ws.Range("D34").Value = "'0.07718"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.19%  "

$ws.Range("D35").Value = "'8.381"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.62%  "

$ws.Range("D36").Value = "'10.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.64%  "

$ws.Range("D37").Value = "'0.05738"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.58%  "

$ws.Range("D38").Value = "'1.001"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.26%  "

$ws.Range("D39").Value = "'4.743"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.76%  "

$ws.Range("D40").Value = "'1.135"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.82%  "

$ws.Range("D41").Value = "'0.1924"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.36%  "

$ws.Range("D42").Value = "'0.02033"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.62%  "

$ws.Range("D43").Value = "'1.337"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.66%  "

$ws.Range("D44").Value = "'3.589"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.00%  "

$ws.Range("D45").Value = "'0.5350"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.76%  "

$ws.Range("D46").Value = "'12.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.05%  "

$ws.Range("D47").Value = "'0.5188"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.21%  "

$ws.Range("D48").Value = "'112.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.33%  "

$ws.Range("D49").Value = "'1.790"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.39%  "

$ws.Range("D50").Value = "'1.060"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.97%  "

$ws.Range("E51").Value = "  -0.26%  "
